$d = $word.ActiveDocument

$replacements = @(
    @("2024-05-01 Wednesday", "2024-05-02 Thursday"),
    @("40×37=1480", "13×32=416"),
    @("88×42=3696", "17×66=1122"),
    @("97×88=8536", "62×87=5394"),
    @("89×28=2492", "44×94=4136"),
    @("88×91=8008", "28×46=1288"),
    @("42×95=3990", "19×22=418"),
    @("48×89=4272", "27×34=918"),
    @("70×66=4620", "52×38=1976"),
    @("18×95=1710", "16×56=896"),
    @("31×78=2418", "94×44=4136"),
    @("52×32=1664", "63×23=1449"),
    @("52×19=988", "62×73=4526"),
    @("28×92=2576", "48×30=1440"),
    @("49×58=2842", "92×69=6348"),
    @("33×74=2442", "66×55=3630"),
    @("37×63=2331", "89×18=1602"),
    @("46×81=3726", "65×23=1495"),
    @("52×82=4264", "86×42=3612"),
    @("38×95=3610", "61×59=3599"),
    @("88×44=3872", "61×91=5551"),
    @("94×41=3854", "71×23=1633"),
    @("27×14=378", "37×55=2035"),
    @("66×22=1452", "88×56=4928"),
    @("29×71=2059", "25×36=900"),
    @("74×13=962", "47×90=4230")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
